# The workbook tracked "Year of Treatment" as its own column (old column B),
# but that information is dropped from this table: the column is removed and
# every column to its right (old C:K) shifts one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Delete()

# The remaining category headers in row 1 (now B1:J1, formerly C1:K1) get a
# ".jamais.jamais" suffix appended to their titles.
$headerRange = $ws.Range("B1:J1")
for ($i = 1; $i -le 9; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $cell.Value2 = $cell.Value2 + ".jamais.jamais"
}
